$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of perft results for version v3 (magic bitboards)
$ws.Range("A4").Value = "v3"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.08
$ws.Range("F4").Value = 1.25
$ws.Range("G4").Value = 29.29
$ws.Range("H4").Value = 30.68

# Update selection to match the new last cell, mirroring the diff's selection change
$ws.Range("H4").Select()
